# Reproduce the recorded edit:
#  1. Slide 16's table switches from the custom "Table_0" style to the
#     built-in table style {0F355CF1-0A48-4CA2-86B8-687AE83A77C4}.
#  2. The presentation's theme colours (the "Integral" palette used by the
#     slide master) are swapped for the plain "Office Theme" palette that
#     used to live only on the notes master.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{0F355CF1-0A48-4CA2-86B8-687AE83A77C4}")

# --- 2. Swap theme colour scheme (Integral -> Office Theme) ---------------
function Convert-HexToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToVbaRgb $officeThemeColors[$i - 1]
}
